$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: row 5 (f34f748b... handback entry) datetimes move forward
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-19 07:20:46"
$wsZhCn.Range("G5").Value = "2016-01-19 07:21:27"

# "de-de" sheet: row 5 (f34f748b... handback entry) datetimes move forward
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-19 07:20:56"
$wsDeDe.Range("G5").Value = "2016-01-19 07:21:44"
